$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PUBGY")

# Helper: writes a contiguous horizontal block of values (columns D..J, i.e. the
# 7 fiscal-year columns) into one row in a single COM call.
function Set-RowValues {
    param($ws, $row, $startCol, $values)
    $n = $values.Count
    $arr = New-Object "object[,]" 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $startColIndex = [int][char]$startCol - [int][char]"A" + 1
    $endColIndex = $startColIndex + $n - 1
    $startCell = $ws.Cells.Item($row, $startColIndex)
    $endCell = $ws.Cells.Item($row, $endColIndex)
    $ws.Range($startCell, $endCell).Value2 = $arr
}

# --- Income Statement ---
Set-RowValues $ws 8 "D" @(11495900, 10920300, 10772200, 8140000, 7801200, 7416400, 6525500)
Set-RowValues $ws 14 "D" @(256900, 1612300, 30300, 79700, 3400, 23600, -2200)
Set-RowValues $ws 15 "D" @(443200, 274900, 295100, 197500, 189600, 333200, 273800)
Set-RowValues $ws 17 "D" @(10019400, 10910200, 9226100, 6940600, 6541200, 6243900, 5502200)
Set-RowValues $ws 18 "D" @(1476500, 10100, 1546100, 1199400, 1260000, 1172500, 1023300)
Set-RowValues $ws 20 "D" @(-29200, -84100, 22400, 22400, 79700, 123400, 34800)
Set-RowValues $ws 21 "D" @(1839900, 201500, 1864300, 1500800, 1534200, 1500600, "NA")
Set-RowValues $ws 22 "D" @(113300, 120100, 122300, 53900, 103200, 159300, 99900)
Set-RowValues $ws 23 "D" @(1334000, -194100, 1446200, 1168000, 1236400, 1136600, 958200)
Set-RowValues $ws 24 "D" @(350100, 383700, 433100, 348900, 334400, 313000, 276000)
Set-RowValues $ws 26 "D" @(984000, -577800, 1013200, 819100, 902100, 823500, 682200)
Set-RowValues $ws 27 "D" @(967200, -591300, 1010900, 807800, 888600, 821300, 668700)
Set-RowValues $ws 32 "D" @(29200, 84100, -22400, -22400, -79700, -123400, -34800)
Set-RowValues $ws 33 "D" @(967200, -591300, 1010900, 807800, 888600, 821300, 668700)
Set-RowValues $ws 35 "D" @(967200, -591300, 1010900, 807800, 888600, 821300, 668700)

# --- Balance Sheet ---
Set-RowValues $ws 41 "D" @(5401300, 2499800, 1876000, 3543200, 3235800, 2948600, 4878400)
Set-RowValues $ws 43 "D" @(23335100, 12014300, 11783100, 9280000, 16854500, 16677300, 15723600)
Set-RowValues $ws 44 "D" @(432000, 455500, 461100, 359000, 344500, 767400, 769700)

# Row 45 ("Other Current Assets") keeps "NA" in columns E, F, G and J; only D, H, I change.
$ws.Range("D45").Value2 = 69600
$ws.Range("H45").Value2 = 255800
$ws.Range("I45").Value2 = 270400
Set-RowValues $ws 46 "D" @(14869700, 14969600, 14120200, 13182300, 10389600, 10196600, 10685800)
Set-RowValues $ws 47 "D" @(362400, 301800, 325400, 259200, 167200, 568800, 301800)
Set-RowValues $ws 48 "D" @(1323900, 718100, 740500, 619300, 1151200, 1135500, 556500)
Set-RowValues $ws 49 "D" @(12003000, 11775300, 13185600, 8932200, 8977000, 14920200, 13894700)
Set-RowValues $ws 52 "D" @(145900, 168300, 178400, 149200, 141400, 107700, 92000)
Set-RowValues $ws 54 "D" @(26680900, 27933100, 28550200, 23142200, 19197200, 18629500, 18456700)
Set-RowValues $ws 57 "D" @(12948900, 13454900, 13201300, 10816000, 9689500, 9255300, 8689800)
Set-RowValues $ws 58 "D" @(444300, 70700, 81900, 397200, 536300, 968300, 1880500)
Set-RowValues $ws 59 "D" @(4778600, 2853200, 3088800, 2342700, 4075100, 4284900, 3699200)
Set-RowValues $ws 60 "D" @(15743800, 16378800, 16372100, 13555900, 12090600, 11795500, 11593500)
Set-RowValues $ws 61 "D" @(3119100, 3220100, 3212300, 1571900, 397200, 562100, 1638100)
Set-RowValues $ws 62 "D" @(1796300, 1529300, 1579800, 1153400, 1363200, 1565200, 1357600)
Set-RowValues $ws 66 "D" @(19998300, 21139400, 21194400, 16313700, 13480700, 13451500, 14081000)
Set-RowValues $ws 72 "D" @(9407900, 2559300, 3433300, 3095600, 8195000, 1993800, 1551700)
Set-RowValues $ws 76 "D" @(6682600, 6793600, 7355800, 6828400, 5716500, 5178000, 4375800)

# --- Cash Flow Statement ---
Set-RowValues $ws 81 "D" @(967200, -591300, 1010900, 807800, 888600, 821300, 668700)
Set-RowValues $ws 83 "D" @(391600, 274900, 295100, 278300, 194100, 204200, "NA")
Set-RowValues $ws 89 "D" @(1668400, 1202800, 1576400, 1159000, 1547200, 1157900, 997400)
Set-RowValues $ws 91 "D" @(-152600, -194100, -259200, -151500, -141400, -138000, -130200)
Set-RowValues $ws 94 "D" @(-467900, -461100, -3937100, -656400, -844900, -683300, "NA")
Set-RowValues $ws 96 "D" @(-190700, -216500, -269300, -124500, -121200, -133500, -144700)
Set-RowValues $ws 100 "D" @(-576700, -265900, 511600, 1311600, -454400, -1431700, "NA")
Set-RowValues $ws 101 "D" @(-425200, 141400, 189600, 149200, -140200, -7900, "NA")
Set-RowValues $ws 102 "D" @(198600, 617100, -1659400, 1963500, 107700, -964900, 20200)
